# Admin Interests Dashboard redesign: add Event Price Amount / Event Price
# Currency columns, and append 4 new interest rows submitted by shiv/sawale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns K1:L1 (same bold/centered/bordered header style as A1:J1) ---
$ws.Range("K1").Value = "Event Price Amount"
$ws.Range("L1").Value = "Event Price Currency"
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)

# --- New data rows (14-17) ---
$newRows = @(
    @("shiv", "sawale", "sawashi01", "Shivam.Sawale@bakerhughes.com", "Intern", "hyderabad", "India", "The Energy Summit / Future Energy Meet 2026", "TBD", "2026-01-06 14:14:34"),
    @("shiv", "sawale", "sawashi01", "Shivam.Sawale@bakerhughes.com", "Intern", "hyderabad", "India", "World Conf. on Robotics and Automation", "TBD", "2026-01-06 14:18:46"),
    @("shiv", "sawale", "sawashi01", "Shivam.Sawale@bakerhughes.com", "Intern", "hyderabad", "India", "The 11th Asia-Pacific Optical Sensors Conference (APOS),", "TBD", "2026-01-06 14:22:33"),
    @("shiv", "sawale", "sawashi01", "Shivam.Sawale@bakerhughes.com", "Intern", "hyderabad", "India", "Gartner Security & Risk Management Summit 2026 (US)", "TBD", "2026-01-06 14:42:04")
)

$rowIndex = 14
foreach ($row in $newRows) {
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}

# --- Empty (but present/typed) string cells for K/L across all data rows (2-17) ---
# A bare empty string is normalized away to a truly blank cell by the engine,
# so force text-type via a quote-prefixed entry, then strip the quote-prefix
# style back off (leaves a real empty-string cell with no extra formatting).
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 11).Value = "'"
    $ws.Cells.Item($r, 11).Style = "Normal"
    $ws.Cells.Item($r, 12).Value = "'"
    $ws.Cells.Item($r, 12).Style = "Normal"
}

# --- Refresh the sheet's used-range dimension ---
$ws.Range("A1:L17").Select()
